$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Rename header row: "<Spalte>_old" -> "<Spalte>_FV2310"
#                        "<Spalte>_new" -> "<Spalte>_FV2404"
# ---------------------------------------------------------------------
$headers = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310",
    "diff",
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# ---------------------------------------------------------------------
# 2. Freeze the header row (split below row 1, top-left cell A2)
# ---------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# ---------------------------------------------------------------------
# 3. Turn the used range into an Excel Table ("Table1")
# ---------------------------------------------------------------------
$usedRange = $ws.Range("A1:U55")
$tbl = $ws.ListObjects.Add(1, $usedRange, [System.Reflection.Missing]::Value, 1)
$tbl.Name = "Table1"
